$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 (shifts existing rows 52:85 down to 53:86,
# and extends the used range to A1:T86, matching the dimension change in the diff).
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C52").Value = 'Ñuble'
$ws.Range("D52").Value = 44879
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 'Fruta'
$ws.Range("G52").Value = 100108
$ws.Range("H52").Value = 'Tropicales y subtropicales'
$ws.Range("I52").Value = 100108002
$ws.Range("J52").Value = 'Mango'
$ws.Range("K52").Value = 'Sin especificar'
$ws.Range("L52").Value = 'Primera'
$ws.Range("M52").Value = 60
$ws.Range("N52").Value = 8000
$ws.Range("O52").Value = 8500
$ws.Range("P52").Value = 8250
$ws.Range("Q52").Value = '$/bandeja 4 kilos'
$ws.Range("R52").Value = 'Brasil'
$ws.Range("S52").Value = 2062
$ws.Range("T52").Value = 4
